$wb = $excel.ActiveWorkbook

$newFile = "09f8bea7-a83d-494d-b118-1d41b940bf20.md"
$newStatus = "Handoff transform failed"
$configName = ".localization-config"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fa7ee3fb788259ebc0168834f10f02ca97e7d82c/.localization-config"
$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/71e178fd45324d652c4ce46c9032d51890e9dda8/e2e/"
$zeroDate = "0001-01-01 00:00:00"
$ignored = "Ignored"

# ---------------------------------------------------------------
# Sheet "Overview": just a text/status refresh, hyperlink target
# filename changes too.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFile
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($mdUrlBase + $newFile), "", "", $newFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", $configName) | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn": handoff failed, drop the handback artifact (col C),
# reset the handoff/handback timestamps, flip Include -> Ignored.
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFile
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Clear()
$wsZh.Range("D2").Value = $zeroDate
$wsZh.Range("G2").Value = $zeroDate
$wsZh.Range("H2").Value = $ignored
$wsZh.Range("D3").Value = $zeroDate
$wsZh.Range("G3").Value = $zeroDate
$wsZh.Range("H3").Value = $ignored

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($mdUrlBase + $newFile), "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configUrl, "", "", $configName) | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de": same treatment.
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFile
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Clear()
$wsDe.Range("D2").Value = $zeroDate
$wsDe.Range("G2").Value = $zeroDate
$wsDe.Range("H2").Value = $ignored
$wsDe.Range("D3").Value = $zeroDate
$wsDe.Range("G3").Value = $zeroDate
$wsDe.Range("H3").Value = $ignored

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($mdUrlBase + $newFile), "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configUrl, "", "", $configName) | Out-Null
